# Generate Report for Handoff
# The file "123134c4-c9c3-405d-b63b-17d350e33907.md" has finished translation
# and is now ready to be handed off. Update its status / priority / timestamp
# on every sheet of the localization-status report.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-01 12:14:33"
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

# ---- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-09-01 12:14:29"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33

# ---- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-09-01 12:14:33"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
